$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Update Schedule sheet (recomputed Cost and Unit Cost values for rows 3-5)
$wsSchedule.Range("E3").Value = -10.54730624999998
$wsSchedule.Range("F3").Value = -0.3100325176366838
$wsSchedule.Range("E4").Value = 572.3830124999999
$wsSchedule.Range("F4").Value = 30.28481547619048
$wsSchedule.Range("E5").Value = -168.28144125
$wsSchedule.Range("F5").Value = -4.946544422398589

# Update Detailed sheet (updated forecast/historical Price values, and one Type relabel)
$wsDetailed.Range("B32").Value = -7.85287
$wsDetailed.Range("B33").Value = 13.91991
$wsDetailed.Range("B34").Value = 24.7014
$wsDetailed.Range("C34").Value = "historical"
$wsDetailed.Range("B35").Value = 26.41646
$wsDetailed.Range("B36").Value = 54.0145
$wsDetailed.Range("B37").Value = 50.52352
$wsDetailed.Range("B38").Value = 58.57113
$wsDetailed.Range("B39").Value = 64.31851
$wsDetailed.Range("B40").Value = 67.69738
$wsDetailed.Range("B41").Value = 76.61194
$wsDetailed.Range("B42").Value = 77.94
$wsDetailed.Range("B43").Value = 72.82266
$wsDetailed.Range("B44").Value = 71.72403
$wsDetailed.Range("B45").Value = 63.77459
$wsDetailed.Range("B47").Value = 63.57817
$wsDetailed.Range("B48").Value = 62.04848
$wsDetailed.Range("B49").Value = 61.77267
$wsDetailed.Range("B50").Value = 57.06003
$wsDetailed.Range("B59").Value = 72.00874
$wsDetailed.Range("B60").Value = 77.71132
$wsDetailed.Range("B62").Value = 81.40685999999999
$wsDetailed.Range("B63").Value = 73.20005
$wsDetailed.Range("B64").Value = 53.27499
$wsDetailed.Range("B65").Value = 0.9831800000000001
$wsDetailed.Range("B66").Value = 0.009650000000000001
$wsDetailed.Range("B68").Value = -6.08274
$wsDetailed.Range("B69").Value = -6.96614
$wsDetailed.Range("B70").Value = -8.786479999999999
$wsDetailed.Range("B71").Value = -9.5
$wsDetailed.Range("B72").Value = -15.04586
$wsDetailed.Range("B73").Value = -18.19399
$wsDetailed.Range("B74").Value = -16.79318
$wsDetailed.Range("B75").Value = -22.40926
$wsDetailed.Range("B77").Value = -23.5
$wsDetailed.Range("B78").Value = -21.03239
$wsDetailed.Range("B79").Value = -23.5
$wsDetailed.Range("B80").Value = -18.71477
$wsDetailed.Range("B82").Value = -5.72703
$wsDetailed.Range("B85").Value = 46.59793
$wsDetailed.Range("B87").Value = 69.05019
$wsDetailed.Range("B88").Value = 100.01
$wsDetailed.Range("B90").Value = 87.00089
$wsDetailed.Range("B91").Value = 73.78428
$wsDetailed.Range("B92").Value = 70.50577
$wsDetailed.Range("B94").Value = 60.25514
